$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: Collapse the three detailed "CORE COMPETENCIES" paragraphs into a
# single summary paragraph that just lists the three category headers.
# ---------------------------------------------------------------------------
$bullet = [char]0x2022
$summary = "Product Management & Strategy " + $bullet + " Technical Product Development " + $bullet + " Platform & Infrastructure"

# Find the "CORE COMPETENCIES" heading paragraph, then operate on the three
# paragraphs that immediately follow it.
$coreHeadingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7) -eq "CORE COMPETENCIES") {
        $coreHeadingIndex = $i
        break
    }
}

$firstDetailIndex = $coreHeadingIndex + 1

# Replace the first of the three paragraphs with the condensed summary text.
$d.Paragraphs.Item($firstDetailIndex).Range.Text = $summary

# The two following paragraphs (Technical Product Development / Platform &
# Infrastructure details) are no longer needed; delete them outright. After
# each delete the following paragraph shifts up into the same index.
$d.Paragraphs.Item($firstDetailIndex + 1).Range.Delete()
$d.Paragraphs.Item($firstDetailIndex + 1).Range.Delete()

# ---------------------------------------------------------------------------
# Change 2: Append a new "TECHNICAL SKILLS" section at the end of the
# document, with a Heading2 title and three condensed-detail paragraphs.
# ---------------------------------------------------------------------------
$lastParaIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastParaIndex)
$lastPara.Range.InsertParagraphAfter()

$headingPara = $d.Paragraphs.Item($lastParaIndex + 1)
$headingPara.Range.Text = "TECHNICAL SKILLS"
$headingPara.Style = $d.Styles.Item("Heading2")

$headingPara.Range.InsertParagraphAfter()
$p1 = $d.Paragraphs.Item($lastParaIndex + 2)
$p1.Style = $d.Styles.Item("Normal")
$p1.Range.Text = "PRODUCT MANAGEMENT & STRATEGY Product Conception & Ideation; Product Architecture & Design; Product Lifecycle Management; B2B SaaS Development; Product Strategy; Stakeholder Management; Product Analytics"

$p1.Range.InsertParagraphAfter()
$p2 = $d.Paragraphs.Item($lastParaIndex + 3)
$p2.Style = $d.Styles.Item("Normal")
$p2.Range.Text = "TECHNICAL PRODUCT DEVELOPMENT Full-Stack Development; Cloud Platforms; Big Data Technologies; Database Design; API Development; DevOps & Deployment; System Integration"

$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs.Item($lastParaIndex + 4)
$p3.Style = $d.Styles.Item("Normal")
$p3.Range.Text = "PLATFORM & INFRASTRUCTURE Multi-tenant Architecture; Data Warehousing; Geospatial Platforms; Real-time Systems; Security & Compliance; Monitoring & Analytics; Documentation & Training"
